{"js": "// Locate the first table (General Format Definition table) and its rows.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// The row that currently reads \"17 - Number of explored locations\" is the\n// last row of the table (index 17, since index 0 is the header row).\nconst row17 = rows.items[rows.items.length - 1];\nconst row17Cells = row17.cells;\nrow17Cells.load(\"items\");\nawait context.sync();\n\nconst col1Cell = row17Cells.items[0];\nconst col2Cell = row17Cells.items[1];\nconst col3Cell = row17Cells.items[2];\n\n// --- Update row 17: split the combined \"17 - Number of explored locations\"\n// entry into a plain \"17\" value, with the description/value columns\n// rewritten to describe the *explored* space count, and the bookmark\n// content (previously in this row) removed - it is relocated to the new\n// \"20\" row below.\ncol1Cell.body.clear();\ncol1Cell.body.insertText(\"17\", \"Start\");\n\ncol2Cell.body.clear();\ncol2Cell.body.insertText(\"# of explored spaces\", \"Start\");\n\ncol3Cell.body.clear();\ncol3Cell.body.insertText(\"Int <= width * height\", \"Start\");\n\nawait context.sync();\n\n// --- Append the three new rows (18, 19, 20) after row 17.\ntable.addRows(\"End\", 3, [\n  [\"\", \"\", \"\"],\n  [\"19\", \"# of known spaces\", \"Int < width * height\"],\n  [\"20 \u2013 number of known spaces\", \"Record of the known spaces\", \"\"],\n]);\nawait context.sync();\n\nrows.load(\"items\");\nawait context.sync();\n\nconst newRowCount = rows.items.length;\nconst row18 = rows.items[newRowCount - 3];\nconst row19 = rows.items[newRowCount - 2];\nconst row20 = rows.items[newRowCount - 1];\n\n// Row 18 needs its text split across multiple runs per the source edit:\n//   col1 -> \"1\" + \"8\" + \" - Number of explored locations\"\n//   col2 -> \"Record \" + \"of \" + \"the explored spaces\"\n//   col3 -> \"Valid X,Y location in space \" (single run, trailing space)\nconst row18Cells = row18.cells;\nrow18Cells.load(\"items\");\nawait context.sync();\n\nconst r18c1 = row18Cells.items[0];\nconst r18c2 = row18Cells.items[1];\nconst r18c3 = row18Cells.items[2];\n\nconst ooxmlWrap = (paraXml) =>\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' + paraXml + '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\n\nr18c1.body.insertOoxml(\n  ooxmlWrap(\n    '<w:p>' +\n      '<w:r><w:t>1</w:t></w:r>' +\n      '<w:r><w:t>8</w:t></w:r>' +\n      '<w:r><w:t xml:space=\"preserve\"> \\u2013 Number of explored locations</w:t></w:r>' +\n    '</w:p>'\n  ),\n  \"Replace\"\n);\n\nr18c2.body.insertOoxml(\n  ooxmlWrap(\n    '<w:p>' +\n      '<w:r><w:t xml:space=\"preserve\">Record </w:t></w:r>' +\n      '<w:r><w:t xml:space=\"preserve\">of </w:t></w:r>' +\n      '<w:r><w:t>the explored spaces</w:t></w:r>' +\n    '</w:p>'\n  ),\n  \"Replace\"\n);\n\nr18c3.body.insertOoxml(\n  ooxmlWrap(\n    '<w:p>' +\n      '<w:r><w:t xml:space=\"preserve\">Valid X,Y location in space </w:t></w:r>' +\n    '</w:p>'\n  ),\n  \"Replace\"\n);\n\nawait context.sync();\n\n// Row 20's \"Valid Values\" column gets the relocated bookmark, with text\n// \"Valid X,Y location in space\" (no trailing space run this time).\nconst row20Cells = row20.cells;\nrow20Cells.load(\"items\");\nawait context.sync();\n\nconst r20c3 = row20Cells.items[2];\nr20c3.body.insertOoxml(\n  ooxmlWrap(\n    '<w:p>' +\n      '<w:r><w:t>Valid X,Y location in space</w:t></w:r>' +\n      '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n      '<w:bookmarkEnd w:id=\"0\"/>' +\n    '</w:p>'\n  ),\n  \"Replace\"\n);\n\nawait context.sync();\n", "ps1": "# Locate the first table (General Format Definition table).\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Helper used to set a table cell's contents from an OOXML paragraph\n# fragment. Using InsertXML (instead of Range.Text) fully replaces the\n# cell's existing content - including any bookmarks or extra runs -\n# and lets us control run boundaries precisely (needed for the rows\n# below that split their text across multiple <w:r> elements).\nfunction Set-CellXml {\n    param($cell, $paraXml)\n    $wrapped = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body>' + $paraXml + '</w:body>' +\n        '</w:document>' +\n        '</pkg:xmlData></pkg:part></pkg:package>'\n    $cell.Range.InsertXML($wrapped)\n}\n\n# --- Row 17 (the last existing row, \"17 - Number of explored locations\")\n# is rewritten to plain \"17\", with the description/valid-value columns\n# now describing the *explored* space count. The bookmark that used to\n# live in this row's third column is relocated to the new \"20\" row\n# created below.\n$row17 = $t.Rows.Last\n$row17Index = $row17.Index\n\n$cell = $t.Cell($row17Index, 1)\n$xml = '<w:p><w:r><w:t>17</w:t></w:r></w:p>'\nSet-CellXml $cell $xml\n\n$cell = $t.Cell($row17Index, 2)\n$xml = '<w:p><w:r><w:t># of explored spaces</w:t></w:r></w:p>'\nSet-CellXml $cell $xml\n\n$cell = $t.Cell($row17Index, 3)\n$xml = '<w:p><w:r><w:t>Int &lt;= width * height</w:t></w:r></w:p>'\nSet-CellXml $cell $xml\n\n# --- Append three new rows (18, 19, 20) after row 17.\n$t.Rows.Add() | Out-Null\n$t.Rows.Add() | Out-Null\n$t.Rows.Add() | Out-Null\n\n$row18Index = $row17Index + 1\n$row19Index = $row17Index + 2\n$row20Index = $row17Index + 3\n\n# Row 18: text split across multiple runs, mirroring the source edit -\n#   col1 -> \"1\" + \"8\" + \" \u2013 Number of explored locations\"\n#   col2 -> \"Record \" + \"of \" + \"the explored spaces\"\n#   col3 -> \"Valid X,Y location in space \" (single run, trailing space)\n$cell = $t.Cell($row18Index, 1)\n$xml = '<w:p>' +\n    '<w:r><w:t>1</w:t></w:r>' +\n    '<w:r><w:t>8</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> \u2013 Number of explored locations</w:t></w:r>' +\n    '</w:p>'\nSet-CellXml $cell $xml\n\n$cell = $t.Cell($row18Index, 2)\n$xml = '<w:p>' +\n    '<w:r><w:t xml:space=\"preserve\">Record </w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\">of </w:t></w:r>' +\n    '<w:r><w:t>the explored spaces</w:t></w:r>' +\n    '</w:p>'\nSet-CellXml $cell $xml\n\n$cell = $t.Cell($row18Index, 3)\n$xml = '<w:p><w:r><w:t xml:space=\"preserve\">Valid X,Y location in space </w:t></w:r></w:p>'\nSet-CellXml $cell $xml\n\n# Row 19: plain new row describing the known-space count.\n$cell = $t.Cell($row19Index, 1)\n$xml = '<w:p><w:r><w:t>19</w:t></w:r></w:p>'\nSet-CellXml $cell $xml\n\n$cell = $t.Cell($row19Index, 2)\n$xml = '<w:p><w:r><w:t># of known spaces</w:t></w:r></w:p>'\nSet-CellXml $cell $xml\n\n$cell = $t.Cell($row19Index, 3)\n$xml = '<w:p><w:r><w:t>Int &lt; width * height</w:t></w:r></w:p>'\nSet-CellXml $cell $xml\n\n# Row 20: combined \"20 - number of known spaces\" row. Its Valid Values\n# column receives the bookmark that was removed from row 17 above (text\n# has no trailing space run this time).\n$cell = $t.Cell($row20Index, 1)\n$xml = '<w:p><w:r><w:t>20 \u2013 number of known spaces</w:t></w:r></w:p>'\nSet-CellXml $cell $xml\n\n$cell = $t.Cell($row20Index, 2)\n$xml = '<w:p><w:r><w:t>Record of the known spaces</w:t></w:r></w:p>'\nSet-CellXml $cell $xml\n\n$cell = $t.Cell($row20Index, 3)\n$xml = '<w:p>' +\n    '<w:r><w:t>Valid X,Y location in space</w:t></w:r>' +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n    '<w:bookmarkEnd w:id=\"0\"/>' +\n    '</w:p>'\nSet-CellXml $cell $xml\n"}
